# Add new columns I (I0) and J (IF) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, styled like the existing header row (bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for rows 2..40, columns I and J.
$IVals = @(9,5,6,1,6,1,8,7,6,6,8,8,9,6,9,8,6,5,8,8,9,9,5,8,8,6,7,9,6,8,7,7,6,8,7,5,8,7,6)
$JVals = @(9,5,6,1,7,1,8,7,6,7,9,8,9,7,9,8,7,6,9,8,9,9,5,8,8,6,7,9,7,8,8,7,7,8,8,5,8,7,6)

for ($i = 0; $i -lt $IVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $IVals[$i]
    $ws.Cells.Item($row, 10).Value = $JVals[$i]
}
